$d = $word.ActiveDocument

$replacements = @(
    @("622÷7=88, 6", "710÷4=177, 2"),
    @("577÷5=115, 2", "460÷8=57, 4"),
    @("300÷5=60, 0", "494÷6=82, 2"),
    @("255÷7=36, 3", "492÷3=164, 0"),
    @("237÷2=118, 1", "283÷9=31, 4"),
    @("786÷7=112, 2", "998÷4=249, 2"),
    @("599÷8=74, 7", "745÷5=149, 0"),
    @("159÷9=17, 6", "363÷3=121, 0"),
    @("311÷8=38, 7", "608÷2=304, 0"),
    @("497÷9=55, 2", "119÷4=29, 3"),
    @("852÷3=284, 0", "861÷7=123, 0"),
    @("696÷9=77, 3", "550÷4=137, 2"),
    @("688÷3=229, 1", "260÷7=37, 1"),
    @("454÷7=64, 6", "213÷9=23, 6"),
    @("978÷8=122, 2", "988÷6=164, 4"),
    @("787÷8=98, 3", "896÷2=448, 0"),
    @("434÷9=48, 2", "352÷5=70, 2"),
    @("431÷9=47, 8", "970÷5=194, 0"),
    @("750÷3=250, 0", "188÷5=37, 3"),
    @("297÷8=37, 1", "768÷2=384, 0"),
    @("159÷7=22, 5", "437÷7=62, 3"),
    @("768÷8=96, 0", "684÷5=136, 4"),
    @("548÷9=60, 8", "898÷4=224, 2"),
    @("174÷2=87, 0", "341÷8=42, 5"),
    @("289÷9=32, 1", "880÷9=97, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
